# Daily attendance processing - reorder "Recorded By" (column G) names so
# that any "System" entry is listed first among the recorders for a session.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -eq $null) { continue }

    $text = [string]$value
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Count -lt 2) { continue }

    if ($parts[0].Trim().ToLower() -eq "system") { continue }

    $tmp = $parts[0]
    $parts[0] = $parts[1]
    $parts[1] = $tmp

    $cell.Value = [string]::Join(", ", $parts)
}
